# Weekly fruit/vegetable price update: insert a new daily record row above
# row 268 (pushing the existing rows 268-283 down to 269-284), then fill
# the newly inserted row with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 268 (existing row 268 and below shift
# down by one row; formatting is inherited from the row above, matching
# the D-column date style already used throughout the table).
$ws.Rows.Item(268).Insert()

# Populate the newly inserted row 268 with the new record.
$ws.Range("A268").Value = 6
$ws.Range("B268").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C268").Value = "Metropolitana"
$ws.Range("D268").Value = 44516
$ws.Range("E268").Value = 13
$ws.Range("F268").Value = 100112030
$ws.Range("G268").Value = "Poroto granado"
$ws.Range("H268").Value = "Sin especificar"
$ws.Range("I268").Value = "Primera"
$ws.Range("J268").Value = 200
$ws.Range("K268").Value = 30000
$ws.Range("L268").Value = 35000
$ws.Range("M268").Value = 32000
$ws.Range("N268").Value = "`$/malla 25 kilos"
$ws.Range("O268").Value = "Perú"
$ws.Range("P268").Value = 1280
$ws.Range("Q268").Value = 25
$ws.Range("R268").Value = "Hortaliza"
